$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 down to rows 3-11 (columns B:G),
# since row 11's original data is dropped and a brand new row of
# data is inserted at row 2.
for ($r = 10; $r -ge 2; $r--) {
    $destRow = $r + 1
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# Insert the new values for row 2
$ws.Range("B2").Value2 = 0.1724578193461484
$ws.Range("C2").Value2 = 0.39058239716261
$ws.Range("D2").Value2 = 0.3033305724894426
$ws.Range("E2").Value2 = 0.550754548314803
$ws.Range("F2").Value2 = 0.5414156770869448
$ws.Range("G2").Value2 = 15
